# Update column C ("Förändrad") date serial value from 45177 (2023-09-08)
# to 45178 (2023-09-09) for every data row (rows 2-33) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
